# AfDD_2023_Annex_Table_Sources.xlsx - "Sources" sheet refresh:
# update several "Latest update" dates and one source-year label to
# reflect newer data pulls, per the author's upload commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sources")

# Global Knowledge Partnership on Migration and Development (KNOMAD)
$ws.Range("D10").Value = "June 2023"

# IMF Investment and Capital Stock (ICSD)
$ws.Range("D13").Value = "Updated 15/06/2022"

# IMF World Economic Outlook Database
$ws.Range("D14").Value = "October 2023"

# Official Development Assistance (ODA) reported by OECD DAC
$ws.Range("D19").Value = "Updated 27/10/2023"

# UNCTADStat Online Data Centre, FDI Online Database
$ws.Range("D24").Value = "Updated 22/09/2022"

# World Development Indicators
$ws.Range("D29").Value = "Updated 10/10/2023"

# World Urbanization Prospects source-year label (2022 -> 2018 edition)
$ws.Range("B32").Value = "World Urbanization Prospects 2018"
